$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Authors) values with re-formatted author strings
$ws.Range("E2").Value = '[Raef%Fadel%NULL%0,   Austin R%Morrison%NULL%0,   Amit%Vahia%NULL%0,   Zachary R%Smith%NULL%0,   Zohra%Chaudhry%NULL%0,   Pallavi%Bhargava%NULL%0,   Joseph%Miller%NULL%0,   Rachel M%Kenney%NULL%0,   George%Alangaden%NULL%0,   Mayur S%Ramesh%mramesh1@hfhs.org%0,   Varidhi%Nauriyal%NULL%0,   Jayanth%Lakshmikanth%NULL%0,   Asif%Abdul Hamed%NULL%0,   Owais%Nadeem%NULL%0,   Kristin%Griebe%NULL%0,   Joseph M%Johnson%NULL%0,   Patrick%Bradley%NULL%0,   Junior%Uduman%NULL%0,   Sara%Hegab%NULL%0,   Jennifer%Swiderek%NULL%0,   Amanda%Godfrey%NULL%0,   Jeffrey%Jennings%NULL%0,   Jayna%Gardner-Gray%NULL%0,   Adam%Ackerman%NULL%0,   Jonathan%Lezotte%NULL%0,   Joseph%Ruhala%NULL%0,   Linoj%Samuel%NULL%0,   Robert J%Tibbetts%NULL%0,   Indira%Brar%NULL%0,   John%McKinnon%NULL%0,   Geehan%Suleyman%NULL%0,   Nicholas%Yared%NULL%0,   Erica%Herc%NULL%0,   Jonathan%Williams%NULL%0,   Odaliz Abreu%Lanfranco%NULL%0,   Anne%Chen%NULL%0,   Marcus%Zervos%NULL%0,   Eric%Scher%NULL%0]'
$ws.Range("E3").Value = '[Xiaowei%Fang%NULL%2,   Qing%Mei%NULL%2,   Tianjun%Yang%NULL%2,   Lei%Li%NULL%2,   Yinzhong%Wang%NULL%2,   Fei%Tong%NULL%2,   Shike%Geng%NULL%2,   Aijun%Pan%NULL%2]'
$ws.Range("E4").Value = '[Xiaowei%Fang%NULL%0,   Qing%Mei%NULL%0,   Tianjun%Yang%NULL%0,   Lei%Li%NULL%0,   Yinzhong%Wang%NULL%0,   Fei%Tong%NULL%0,   Shike%Geng%NULL%0,   Aijun%Pan%NULL%0]'
$ws.Range("E5").Value = '[Yun%Ling%NULL%2,   Shui-Bao%Xu%NULL%2,   Yi-Xiao%Lin%NULL%2,   Di%Tian%NULL%2,   Zhao-Qin%Zhu%NULL%2,   Fa-Hui%Dai%NULL%2,   Fan%Wu%NULL%2,   Zhi-Gang%Song%NULL%2,   Wei%Huang%NULL%2,   Jun%Chen%NULL%0,   Bi-Jie%Hu%NULL%2,   Sheng%Wang%NULL%2,   En-Qiang%Mao%NULL%2,   Lei%Zhu%NULL%0,   Wen-Hong%Zhang%NULL%2,   Hong-Zhou%Lu%NULL%2,   Peng%Lyu%NULL%0,   Peng%Lyu%NULL%0]'
$ws.Range("E6").Value = '[Xiaofan%Lu%NULL%0,   Taige%Chen%NULL%0,   Yang%Wang%NULL%0,   Jun%Wang%NULL%0,   Fangrong%Yan%f.r.yan@163.com%0]'
$ws.Range("E7").Value = '[Yin%Wang%NULL%0,   Weiwei%Jiang%NULL%0,   Qi%He%NULL%0,   Cheng%Wang%NULL%0,   Baoju%Wang%NULL%0,   Pan%Zhou%NULL%0,   Nianguo%Dong%dongnianguo63@gmail.com%0,   Qiaoxia%Tong%2013xh0859@hust.edu.cn%0]'
$ws.Range("E8").Value = '[ J.%Wu%null%1,    J.% Huang%null%1,    G.% Zhu%null%1,    Y.% Liu%null%2,    H.% Xiao%null%1,    Q.% Zhou%null%0,    X.% Si%null%1,    H.% Yi%null%1,    C.% Wang%null%1,    D.% Yang%null%0,    S.% Chen%null%1,    X.% Liu%null%1,    Z.% Liu%null%1,    Q.% Wang%null%1,    Q.% Lv%null%1,    Y.% Huang%null%1,    Y.% Yu%null%1,    X.% Guan%null%1,    Y.% Li%null%1,    K.% Nirantharakumar%null%1,    K.% Cheng%null%1,    S.% Peng%null%1,    H. % Xiao%null%1]'
$ws.Range("E10").Value = '[Xiaobo%Yang%NULL%0,   Yuan%Yu%NULL%0,   Jiqian%Xu%NULL%0,   Huaqing%Shu%NULL%0,   Jia''an%Xia%NULL%0,   Hong%Liu%NULL%0,   Yongran%Wu%NULL%0,   Lu%Zhang%NULL%0,   Zhui%Yu%NULL%0,   Minghao%Fang%NULL%0,   Ting%Yu%NULL%0,   Yaxin%Wang%NULL%0,   Shangwen%Pan%NULL%0,   Xiaojing%Zou%NULL%0,   Shiying%Yuan%NULL%0,   You%Shang%NULL%0]'
$ws.Range("E11").Value = '[Lei%Zha%NULL%0,   Shirong%Li%NULL%0,   Lingling%Pan%NULL%0,   Boris%Tefsen%NULL%0,   Yeshan%Li%NULL%0,   Neil%French%NULL%0,   Liyun%Chen%NULL%0,   Gang%Yang%NULL%0,   Elmer V%Villanueva%Villanueva@xjtlu.edu.cn%0]'
$ws.Range("E12").Value = '[Wei%Zhou%NULL%2,   Yisi%Liu%NULL%2,   Dongdong%Tian%NULL%2,   Cheng%Wang%NULL%0,   Sa%Wang%NULL%2,   Jing%Cheng%NULL%2,   Ming%Hu%huming74@163.com%0,   Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0,   Yue%Gao%gaoyue@bmi.ac.cn%2]'
$ws.Range("E13").Value = '[Wei%Zhou%NULL%0,   Yisi%Liu%NULL%0,   Dongdong%Tian%NULL%0,   Cheng%Wang%NULL%0,   Sa%Wang%NULL%0,   Jing%Cheng%NULL%0,   Ming%Hu%huming74@163.com%0,   Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0,   Yue%Gao%gaoyue@bmi.ac.cn%0]'
$ws.Range("E14").Value = '[Rong-chang%Chen%NULL%1,   Xiao-ping%Tang%NULL%1,   Shou-yong%Tan%NULL%1,   Bi-ling%Liang%NULL%1,   Zhuo-yue%Wan%NULL%1,   Ji-qian%Fang%NULL%1,   Nanshan%Zhong%NULL%1]'
$ws.Range("E15").Value = '[Justin%Coyle%NULL%1,   Efehi%Igbinomwanhia%NULL%1,   Alejandro%Sanchez-Nadales%alejandro.sanchez2@advocatehealth.com%1,   Sorin%Danciu%NULL%1,   Chae%Chu%NULL%1,   Nishit%Shah%NULL%1]'
$ws.Range("E16").Value = '[Yuanyuan%Han%NULL%1,   Mao%Jiang%NULL%1,   Da%Xia%NULL%1,   Lichao%He%NULL%1,   Xin%Lv%NULL%1,   Xiaohua%Liao%NULL%1,   Jie%Meng%mengjie@csu.edu.cn%1]'
$ws.Range("E17").Value = '[Kui%Liu%NULL%0,   Yuan-Yuan%Fang%NULL%0,   Yan%Deng%NULL%0,   Wei%Liu%NULL%0,   Mei-Fang%Wang%NULL%0,   Jing-Ping%Ma%NULL%0,   Wei%Xiao%NULL%0,   Ying-Nan%Wang%NULL%0,   Min-Hua%Zhong%NULL%0,   Cheng-Hong%Li%NULL%0,   Guang-Cai%Li%NULL%0,   Hui-Guo%Liu%NULL%0,   Xiu-Yuan%Hao%NULL%0,   Xiu-Yuan%Hao%NULL%0,   Pei-Fang%Wei%NULL%0]'

# Row 8: F and G columns change (ID / ID Format)
$ws.Range("F8").Value = 'not found'
$ws.Range("G8").Value = 'N/A'

# New column I (Other found locations)
$ws.Range("I1").Value = 'Other found locations'
$ws.Range("I2").Value = '_PMC'
$ws.Range("I3").Value = '_PMC_elsevier'
$ws.Range("I4").Value = '_PMC_elsevier'
$ws.Range("I5").Value = '_PMC'
$ws.Range("I6").Value = '_PMC_Springer'
$ws.Range("I7").Value = '_PMC_Springer'
$ws.Range("I8").Value = '_MedBiorxiv'
$ws.Range("I9").Value = ''
$ws.Range("I10").Value = '_PMC_elsevier'
$ws.Range("I11").Value = '_PMC'
$ws.Range("I12").Value = '_PMC_Springer'
$ws.Range("I13").Value = '_PMC_Springer'
$ws.Range("I14").Value = '_PMC_elsevier'
$ws.Range("I15").Value = '_PMC_elsevier'
$ws.Range("I16").Value = '_PMC_elsevier'
$ws.Range("I17").Value = '_PMC'
$ws.Range("I18").Value = ''
